$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"

# ---------------------------------------------------------------------------
# Paragraph 1: "]Greetings from Etiqa!"  ->  "Greetings from Etiqa!"
#   - drop the stray "]" run
#   - split "Greetings from Etiqa!" into 3 runs, wrapping "Etiqa" with
#     proofErr spellStart/spellEnd (simulating Word's spell-checker)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    if ($full.Text -eq "]Greetings from Etiqa!`r") {
        $xml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="45ACBB6B" w14:textId="6A44830D" w:rsidR="006C66B5" w:rsidRDefault="00882133" w:rsidP="006C66B5"><w:pPr><w:pStyle w:val="p1"/></w:pPr><w:r><w:t xml:space="preserve">Greetings from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Etiqa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>!</w:t></w:r></w:p>'
        $full.InsertXML($xml)

        $p2 = $d.Paragraphs.Item($i)
        $r2 = $p2.Range
        $charRange = $d.Range($r2.Start, $r2.End - 1)
        $charRange.Style = "s1"
        break
    }
}

# ---------------------------------------------------------------------------
# Paragraph 2: "Etiqa Life Insurance (Cambodia) Plc."
#   - split into "Etiqa" + " Life Insurance (Cambodia) Plc." runs, wrapping
#     "Etiqa" with proofErr spellStart/spellEnd
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    if ($full.Text -eq "Etiqa Life Insurance (Cambodia) Plc.`r") {
        $xml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="39C6D915" w14:textId="77777777" w:rsidR="006C66B5" w:rsidRDefault="006C66B5" w:rsidP="006C66B5"><w:pPr><w:pStyle w:val="p1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Etiqa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Life Insurance (Cambodia) Plc.</w:t></w:r></w:p>'
        $full.InsertXML($xml)

        $p2 = $d.Paragraphs.Item($i)
        $r2 = $p2.Range
        $charRange = $d.Range($r2.Start, $r2.End - 1)
        $charRange.Style = "s1"
        break
    }
}
